$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "females are more likely to purchase",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "females are slightly more likely to purchase",
    2
)
